# edit.ps1 - applies the churn-analysis presentation update described by the
# target diff:
#   * slide 18 (picture "annual_financial_overview.png"): off.x 12739133 -> 12739132
#   * slide 21 (textbox "80%"): off.x 1269999 -> 1270000, ext.cx 21844001 -> 21844000,
#       ext.cy 4488605 -> 4488604
#   * slide 22 (body textbox id=233): run text "$8,743.50 in monthly revenue at risk"
#       -> "$8,743 in monthly revenue at risk" (and the shape's auto-generated Name,
#       which mirrors the shape's text truncated to 254 characters, updates to match)
#   * slide 23 (textbox "Target high-risk customers..."): off.x 1269999 -> 1270000,
#       ext.cx 21844001 -> 21844000
#   * slide 24 (textbox "Closing statement"): off.x 1269999 -> 1270000,
#       ext.cx 21844001 -> 21844000
#
# NOTE on precision: Shape.Left/Top/Width/Height are exposed as single precision
# (float32) points, while OOXML stores position/size in EMU (1 pt = 12700 EMU) as
# integers. The runtime converts float32-points -> EMU by truncation, so naively
# writing "emu / 12700.0" can truncate down to emu-1 because of float32 rounding.
# To land on the exact target EMU we use point literals that are already exact
# float32 values chosen so that value*12700 (truncated) reproduces the desired EMU.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 18: picture "annual_financial_overview.png" (id=222) - nudge x offset
# 12739133 -> 12739132 EMU (y, width, height unchanged)
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$pic = $s18.Shapes.Item(2)
if ($pic.Id -ne 222) { throw "Slide 18 shape 2 id mismatch: $($pic.Id)" }
$pic.Left = 1003.081298828125

# ---------------------------------------------------------------------------
# Slide 21: textbox "80%" (id=229)
# off.x 1269999 -> 1270000 ; ext.cx 21844001 -> 21844000 ; ext.cy 4488605 -> 4488604
# ---------------------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$sh229 = $s21.Shapes.Item(1)
if ($sh229.Id -ne 229) { throw "Slide 21 shape 1 id mismatch: $($sh229.Id)" }
$sh229.Left = 100.0
$sh229.Width = 1720.0
$sh229.Height = 353.43341064453125

# ---------------------------------------------------------------------------
# Slide 22: body textbox (id=233) - update the dollar figure text
# "$8,743.50 in monthly revenue at risk" -> "$8,743 in monthly revenue at risk"
# ---------------------------------------------------------------------------
$s22 = $p.Slides.Item(22)
$sh233 = $s22.Shapes.Item(2)
if ($sh233.Id -ne 233) { throw "Slide 22 shape 2 id mismatch: $($sh233.Id)" }
$tr233 = $sh233.TextFrame.TextRange
$run7 = $tr233.Runs(7, 1)
if ($run7.Text -ne '$8,743.50 in monthly revenue at risk') { throw "Slide 22 run 7 text mismatch: $($run7.Text)" }
$run7.Text = '$8,743 in monthly revenue at risk'

# The shape's Name mirrors its full text truncated to 254 characters; this
# runtime does not recompute it automatically on a text edit, so set it
# explicitly to match the new content.
$sh233.Name = 'Our model identified 109 high-risk customers, representing 7.7% of the test group. These customers account for a potential $8,743 in monthly revenue at risk, which scales to approximately $104,922 annually. On average, high-risk customers have monthly ch'

# ---------------------------------------------------------------------------
# Slide 23: textbox "Target high-risk customers..." (id=236)
# off.x 1269999 -> 1270000 ; ext.cx 21844001 -> 21844000 (ext.cy unchanged)
# ---------------------------------------------------------------------------
$s23 = $p.Slides.Item(23)
$sh236 = $s23.Shapes.Item(2)
if ($sh236.Id -ne 236) { throw "Slide 23 shape 2 id mismatch: $($sh236.Id)" }
$sh236.Left = 100.0
$sh236.Width = 1720.0

# ---------------------------------------------------------------------------
# Slide 24: textbox "Closing statement" (id=238)
# off.x 1269999 -> 1270000 ; ext.cx 21844001 -> 21844000 (ext.cy unchanged)
# ---------------------------------------------------------------------------
$s24 = $p.Slides.Item(24)
$sh238 = $s24.Shapes.Item(1)
if ($sh238.Id -ne 238) { throw "Slide 24 shape 1 id mismatch: $($sh238.Id)" }
$sh238.Left = 100.0
$sh238.Width = 1720.0
